# Append new Lancers job listing captured at 2025-10-02 01:41:52 JST,
# and refresh the "取得日時" (fetched-at) timestamp on every existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-02 01:41:52"

# ---------------------------------------------------------------------
# 1) Remove all existing hyperlinks up-front. Inserting a row shifts
#    cell *contents* but this runtime leaves the Hyperlinks collection
#    anchored to the old cell addresses, so the safest way to keep
#    everything consistent is to rebuild the hyperlinks from scratch
#    after the row insert/content updates below.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2) Insert a new row above the current row 17 ("CentOS..." listing),
#    pushing the former rows 17-20 down to 18-21.
# ---------------------------------------------------------------------
$ws.Rows.Item(17).Insert()

# ---------------------------------------------------------------------
# 3) Populate the newly inserted row 17 with the new job listing.
# ---------------------------------------------------------------------
$ws.Range("B17").Value = "【在宅勤務】ランサーズ業務委託で働ける、ネパール人個人の方を募集します!"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5404906"
$ws.Range("G17").Value = 18

# ---------------------------------------------------------------------
# 4) Refresh the timestamp in column A for every data row (2-21), the
#    new row included.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("A$r").Value = $newTimestamp
}

# ---------------------------------------------------------------------
# 5) Rebuild the hyperlinks for column F across every data row so the
#    link targets line up with the (possibly shifted) row contents.
# ---------------------------------------------------------------------
$linkTargets = @{
    2  = "https://www.lancers.jp/work/detail/5391872"
    3  = "https://www.lancers.jp/work/detail/5404026"
    4  = "https://www.lancers.jp/work/detail/5398112"
    5  = "https://www.lancers.jp/work/detail/5404680"
    6  = "https://www.lancers.jp/work/detail/5404305"
    7  = "https://www.lancers.jp/work/detail/5403988"
    8  = "https://www.lancers.jp/work/detail/5404059"
    9  = "https://www.lancers.jp/work/detail/5404342"
    10 = "https://www.lancers.jp/work/detail/5404426"
    11 = "https://www.lancers.jp/work/detail/5371747"
    12 = "https://www.lancers.jp/work/detail/5404650"
    13 = "https://www.lancers.jp/work/detail/5341051"
    14 = "https://www.lancers.jp/work/detail/5367840"
    15 = "https://www.lancers.jp/work/detail/5372984"
    16 = "https://www.lancers.jp/work/detail/5391221"
    17 = "https://www.lancers.jp/work/detail/5404906"
    18 = "https://www.lancers.jp/work/detail/5404155"
    19 = "https://www.lancers.jp/work/detail/5404010"
    20 = "https://www.lancers.jp/work/detail/5404730"
    21 = "https://www.lancers.jp/work/detail/5404652"
}

for ($r = 2; $r -le 21; $r++) {
    $ws.Hyperlinks.Add($ws.Range("F$r"), $linkTargets[$r])
}
